# Sync attendance_reports: reorder "Recorded By" (column G) values so that
# "System" appears first in the comma-separated list of recorders, while
# preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.Length -eq 0) { continue }
    if (-not $val.Contains("System")) { continue }

    $parts = $val -split ",\s*"

    if ($parts.Count -le 1) { continue }
    if ($parts[0].Equals("System")) { continue }

    # find the FIRST exact-case "System" token (not "system" / "SYSTEM" / etc.)
    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $idx = $i
            break
        }
    }

    if ($idx -lt 0) { continue }

    $newVal = "System"
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) {
            $newVal = $newVal + ", " + $parts[$i]
        }
    }

    $cell.Value = $newVal
}
